# RWS project, de omreisroutes incl uitval van het onderliggende net.
# Updates the "input" sheet: renames the pickle-cache file names to the
# "including_underlying" variant, re-orders/renumbers the road ids, and
# appends four new analysis rows (13-16) that re-use the same settings as
# the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# --- 1. Extend rows 13:16 with the same formatting as the existing data rows
# (row 2 is a fully representative template: same style index on every cell).
$ws.Range("A2:Y2").Copy()
$ws.Range("A13:Y16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Re-populate column A (road id) and column P (pickle file name) for
# every data row (2-16). All other columns for rows 2-12 are unchanged; the
# new rows 13-16 re-use the same values as the other rows in those columns.
$ws.Range("A2").Value = 13944
$ws.Range("P2").Value = "including_underlying/road_gdf_sel_incl_underl13944.p"

$ws.Range("A3").Value = 13946
$ws.Range("P3").Value = "including_underlying/road_gdf_sel_incl_underl13946.p"

$ws.Range("A4").Value = 70012
$ws.Range("P4").Value = "including_underlying/road_gdf_sel_incl_underl70012.p"

$ws.Range("A5").Value = 13814
$ws.Range("P5").Value = "including_underlying/road_gdf_sel_incl_underl13814.p"

$ws.Range("A6").Value = 13813
$ws.Range("P6").Value = "including_underlying/road_gdf_sel_incl_underl13813.p"

$ws.Range("A7").Value = 13943
$ws.Range("P7").Value = "including_underlying/road_gdf_sel_incl_underl13943.p"

$ws.Range("A8").Value = 70009
$ws.Range("P8").Value = "including_underlying/road_gdf_sel_incl_underl70009.p"

$ws.Range("A9").Value = 13165
$ws.Range("P9").Value = "including_underlying/road_gdf_sel_incl_underl13165.p"

$ws.Range("A10").Value = 14013
$ws.Range("P10").Value = "including_underlying/road_gdf_sel_incl_underl14013.p"

$ws.Range("A11").Value = 13173
$ws.Range("P11").Value = "including_underlying/road_gdf_sel_incl_underl13173.p"

$ws.Range("A12").Value = 19559
$ws.Range("P12").Value = "including_underlying/road_gdf_sel_incl_underl19559.p"

# --- 3. New rows 13-16: fill in every column, matching the template used by
# the other analysis rows (same analysis settings, only id/pickle differ).
$ws.Range("A13").Value = 13945
$ws.Range("B13").Value = "Redundancy-based criticality"
$ws.Range("C13").Value = "Multi-link Disruption_RWS"
$ws.Range("D13").Value = "Network based on OSM dump"
$ws.Range("F13").Value = "NL_with_margin_from_EU_dump.osm.pbf"
$ws.Range("L13").Value = "drive"
$ws.Range("M13").Value = "motorway"
$ws.Range("N13").Value = "scenario_13794_wgs84.tif"
$ws.Range("O13").Value = "Pavement_avg_depth"
$ws.Range("P13").Value = "including_underlying/road_gdf_sel_incl_underl13945.p"
$ws.Range("R13").Value = "m"
$ws.Range("S13").Value = "mean"
$ws.Range("T13").Value = 0.1

$ws.Range("A14").Value = 13937
$ws.Range("B14").Value = "Redundancy-based criticality"
$ws.Range("C14").Value = "Multi-link Disruption_RWS"
$ws.Range("D14").Value = "Network based on OSM dump"
$ws.Range("F14").Value = "NL_with_margin_from_EU_dump.osm.pbf"
$ws.Range("L14").Value = "drive"
$ws.Range("M14").Value = "motorway"
$ws.Range("N14").Value = "scenario_13794_wgs84.tif"
$ws.Range("O14").Value = "Pavement_avg_depth"
$ws.Range("P14").Value = "including_underlying/road_gdf_sel_incl_underl13937.p"
$ws.Range("R14").Value = "m"
$ws.Range("S14").Value = "mean"
$ws.Range("T14").Value = 0.1

$ws.Range("A15").Value = 13812
$ws.Range("B15").Value = "Redundancy-based criticality"
$ws.Range("C15").Value = "Multi-link Disruption_RWS"
$ws.Range("D15").Value = "Network based on OSM dump"
$ws.Range("F15").Value = "NL_with_margin_from_EU_dump.osm.pbf"
$ws.Range("L15").Value = "drive"
$ws.Range("M15").Value = "motorway"
$ws.Range("N15").Value = "scenario_13794_wgs84.tif"
$ws.Range("O15").Value = "Pavement_avg_depth"
$ws.Range("P15").Value = "including_underlying/road_gdf_sel_incl_underl13812.p"
$ws.Range("R15").Value = "m"
$ws.Range("S15").Value = "mean"
$ws.Range("T15").Value = 0.1

$ws.Range("A16").Value = 19558
$ws.Range("B16").Value = "Redundancy-based criticality"
$ws.Range("C16").Value = "Multi-link Disruption_RWS"
$ws.Range("D16").Value = "Network based on OSM dump"
$ws.Range("F16").Value = "NL_with_margin_from_EU_dump.osm.pbf"
$ws.Range("L16").Value = "drive"
$ws.Range("M16").Value = "motorway"
$ws.Range("N16").Value = "scenario_13794_wgs84.tif"
$ws.Range("O16").Value = "Pavement_avg_depth"
$ws.Range("P16").Value = "including_underlying/road_gdf_sel_incl_underl19558.p"
$ws.Range("R16").Value = "m"
$ws.Range("S16").Value = "mean"
$ws.Range("T16").Value = 0.1

# --- 4. Restore the selected cell shown in the workbook.
$ws.Range("P28").Select()
